$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 02:05"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 1685508
$ws.Range("C4").Value = 18680
$ws.Range("D4").Value = 451702
$ws.Range("E4").Value = 1134520
$ws.Range("G4").Value = 603
$ws.Range("H4").Value = 99286

# Row 5: Brasil -> Brasil
$ws.Range("B5").Value = 363618
$ws.Range("C5").Value = 16220
$ws.Range("E5").Value = 190991
$ws.Range("G5").Value = 703
$ws.Range("H5").Value = 22716

# Row 68: Irak -> Camerun
$ws.Range("A68").Value = "Camerun"
$ws.Range("B68").Value = 4890
$ws.Range("C68").Value = 490
$ws.Range("D68").Value = 1865
$ws.Range("E68").Value = 2860
$ws.Range("G68").Value = 6
$ws.Range("H68").Value = 165

# Row 69: Camerun -> Irak
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 4469
$ws.Range("C69").Value = 197
$ws.Range("D69").Value = 2738
$ws.Range("E69").Value = 1571
$ws.Range("G69").Value = 8
$ws.Range("H69").Value = 160

# Row 144: Ruanda -> Guayana Francesa
$ws.Range("A144").Value = "Guayana Francesa"
$ws.Range("B144").Value = 328
$ws.Range("C144").Value = 49
$ws.Range("D144").Value = 145
$ws.Range("E144").Value = 182
$ws.Range("H144").Value = 1

# Row 145: Vietnam -> Ruanda
$ws.Range("A145").Value = "Ruanda"
$ws.Range("B145").Value = 327
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 237
$ws.Range("E145").Value = 90

# Row 146: Montenegro -> Vietnam
$ws.Range("A146").Value = "Vietnam"
$ws.Range("B146").Value = 325
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 267
$ws.Range("E146").Value = 58
$ws.Range("H146").Value = 0

# Row 147: Guayana Francesa -> Montenegro
$ws.Range("A147").Value = "Montenegro"
$ws.Range("B147").Value = 324
$ws.Range("D147").Value = 315
$ws.Range("E147").Value = 0
$ws.Range("H147").Value = 9

# Row 163: Brunei -> Brunei
$ws.Range("D163").Value = 137
$ws.Range("E163").Value = 3

# Row 181: Zimbabue -> Zimbabue
$ws.Range("D181").Value = 25
$ws.Range("E181").Value = 27

# Row 198: Belice -> Nueva Caledonia
$ws.Range("A198").Value = "Nueva Caledonia"
$ws.Range("D198").Value = 18
$ws.Range("H198").Value = 0

# Row 199: Nueva Caledonia -> Santa Lucia
$ws.Range("A199").Value = "Santa Lucia"

# Row 200: Santa Lucia -> Belice
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

# Row 207: Groenlandia -> Islas Turcas y Caicos
$ws.Range("A207").Value = "Islas Turcas y Caicos"
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 10
$ws.Range("H207").Value = 1

# Row 208: Islas Turcas y Caicos -> Groenlandia
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 0

# Row 210: Montserrat -> Seychelles
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# Row 211: Seychelles -> Montserrat
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 214: Bonaire, San Eustaquio y Saba -> Sahara Occidental
$ws.Range("A214").Value = "Sahara Occidental"

# Row 215: San Bartolome -> Bonaire, San Eustaquio y Saba
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"

# Row 216: Sahara Occidental -> San Bartolome
$ws.Range("A216").Value = "San Bartolome"
